$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plain text-field updates (timestamps, measurements, etc.) - a simple
# Value assignment is enough since Excel does not reinterpret these
# strings as numbers/dates.
$ws.Range("E2").Value = '2026-02-19 18:48:47'
$ws.Range("I2").Value = '2.4 mm'
$ws.Range("E3").Value = '2026-02-19 18:48:50'
$ws.Range("I3").Value = '3.9 mm'
$ws.Range("E4").Value = '2026-02-19 18:48:53'
$ws.Range("J4").Value = '1009.5 hPa'
$ws.Range("L4").Value = '45.4 km/h - 250º 18:08 TU'
$ws.Range("E5").Value = '2026-02-19 18:48:56'
$ws.Range("I5").Value = '7.0 mm'
$ws.Range("E6").Value = '2026-02-19 18:48:58'
$ws.Range("J6").Value = '1009.6 hPa'
$ws.Range("E7").Value = '2026-02-19 18:49:01'
$ws.Range("J7").Value = '1010.5 hPa'
$ws.Range("E8").Value = '2026-02-19 18:49:04'
$ws.Range("J8").Value = '1010.2 hPa'
$ws.Range("K8").Value = '13.3 MJ/m2'
$ws.Range("E9").Value = '2026-02-19 18:49:07'
$ws.Range("E10").Value = '2026-02-19 18:49:09'
$ws.Range("O10").Value = '10.9 °C'
$ws.Range("E11").Value = '2026-02-19 18:49:12'
$ws.Range("O11").Value = '5.4 °C'
$ws.Range("E12").Value = '2026-02-19 18:49:14'
$ws.Range("E13").Value = '2026-02-19 18:49:17'
$ws.Range("J13").Value = '1010.8 hPa'
$ws.Range("E14").Value = '2026-02-19 18:49:20'
$ws.Range("E15").Value = '2026-02-19 18:49:23'
$ws.Range("E16").Value = '2026-02-19 18:49:25'
$ws.Range("I16").Value = '7.7 mm'
$ws.Range("E17").Value = '2026-02-19 18:49:28'
$ws.Range("E18").Value = '2026-02-19 18:49:31'
$ws.Range("J18").Value = '1009.8 hPa'
$ws.Range("O18").Value = '11.9 °C'
$ws.Range("E19").Value = '2026-02-19 18:49:33'
$ws.Range("E20").Value = '2026-02-19 18:49:36'
$ws.Range("E21").Value = '2026-02-19 18:49:39'
$ws.Range("J21").Value = '1010.6 hPa'
$ws.Range("E22").Value = '2026-02-19 18:49:41'
$ws.Range("O22").Value = '-7.7 °C'
$ws.Range("E23").Value = '2026-02-19 18:49:44'
$ws.Range("I23").Value = '8.3 mm'
$ws.Range("E24").Value = '2026-02-19 18:49:47'
$ws.Range("J24").Value = '1014.3 hPa'
$ws.Range("O24").Value = '9.1 °C'
$ws.Range("E25").Value = '2026-02-19 18:49:50'
$ws.Range("I25").Value = '4.0 mm'
$ws.Range("E26").Value = '2026-02-19 18:49:52'
$ws.Range("J26").Value = '1009.5 hPa'
$ws.Range("E27").Value = '2026-02-19 18:49:55'
$ws.Range("E28").Value = '2026-02-19 18:49:58'
$ws.Range("J28").Value = '1009.4 hPa'
$ws.Range("E29").Value = '2026-02-19 18:50:01'
$ws.Range("O29").Value = '11.0 °C'
$ws.Range("E30").Value = '2026-02-19 18:50:04'
$ws.Range("J30").Value = '1009.6 hPa'
$ws.Range("E31").Value = '2026-02-19 18:50:07'
$ws.Range("J31").Value = '1009.0 hPa'
$ws.Range("O31").Value = '11.7 °C'
$ws.Range("E32").Value = '2026-02-19 18:50:09'
$ws.Range("E33").Value = '2026-02-19 18:50:12'
$ws.Range("J33").Value = '1010.4 hPa'
$ws.Range("E34").Value = '2026-02-19 18:50:15'
$ws.Range("E35").Value = '2026-02-19 18:50:18'
$ws.Range("J35").Value = '1015.7 hPa'
$ws.Range("E36").Value = '2026-02-19 18:50:21'
$ws.Range("J36").Value = '1009.9 hPa'
$ws.Range("E37").Value = '2026-02-19 18:50:23'
$ws.Range("J37").Value = '1010.9 hPa'
$ws.Range("O37").Value = '5.8 °C'
$ws.Range("E38").Value = '2026-02-19 18:50:26'
$ws.Range("E39").Value = '2026-02-19 18:50:29'
$ws.Range("I39").Value = '4.4 mm'
$ws.Range("E40").Value = '2026-02-19 18:50:31'
$ws.Range("J40").Value = '1011.9 hPa'
$ws.Range("E41").Value = '2026-02-19 18:50:34'
$ws.Range("J41").Value = '1012.4 hPa'
$ws.Range("E42").Value = '2026-02-19 18:50:37'
$ws.Range("E43").Value = '2026-02-19 18:50:39'
$ws.Range("E44").Value = '2026-02-19 18:50:42'
$ws.Range("I44").Value = '8.2 mm'
$ws.Range("E45").Value = '2026-02-19 18:50:45'
$ws.Range("I45").Value = '3.4 mm'
$ws.Range("J45").Value = '1014.8 hPa'
$ws.Range("E46").Value = '2026-02-19 18:50:48'
$ws.Range("J46").Value = '1015.2 hPa'

# Percentage-like text values (e.g. "67%") need special handling:
# setting .Value = "67%" directly makes Excel COM auto-convert the
# text into the numeric percentage 0.67 with a Percent number format,
# which does not match the source data (plain text cells, General
# format). We briefly force a Text format so the percent string is
# kept literally, then restore the original "General" look by pasting
# the cell format from column F of the same row, which carries the
# same untouched style and is never itself edited.

$ws.Range("H11").NumberFormat = "@"
$ws.Range("H11").Value = '67%'
$ws.Range("F11").Copy()
$ws.Range("H11").PasteSpecial(-4122)

$ws.Range("H13").NumberFormat = "@"
$ws.Range("H13").Value = '64%'
$ws.Range("F13").Copy()
$ws.Range("H13").PasteSpecial(-4122)

$ws.Range("H14").NumberFormat = "@"
$ws.Range("H14").Value = '47%'
$ws.Range("F14").Copy()
$ws.Range("H14").PasteSpecial(-4122)

$ws.Range("H18").NumberFormat = "@"
$ws.Range("H18").Value = '59%'
$ws.Range("F18").Copy()
$ws.Range("H18").PasteSpecial(-4122)

$ws.Range("H20").NumberFormat = "@"
$ws.Range("H20").Value = '87%'
$ws.Range("F20").Copy()
$ws.Range("H20").PasteSpecial(-4122)

$ws.Range("H21").NumberFormat = "@"
$ws.Range("H21").Value = '63%'
$ws.Range("F21").Copy()
$ws.Range("H21").PasteSpecial(-4122)

$ws.Range("H23").NumberFormat = "@"
$ws.Range("H23").Value = '77%'
$ws.Range("F23").Copy()
$ws.Range("H23").PasteSpecial(-4122)

$ws.Range("H25").NumberFormat = "@"
$ws.Range("H25").Value = '66%'
$ws.Range("F25").Copy()
$ws.Range("H25").PasteSpecial(-4122)

$ws.Range("H31").NumberFormat = "@"
$ws.Range("H31").Value = '51%'
$ws.Range("F31").Copy()
$ws.Range("H31").PasteSpecial(-4122)

$ws.Range("H33").NumberFormat = "@"
$ws.Range("H33").Value = '60%'
$ws.Range("F33").Copy()
$ws.Range("H33").PasteSpecial(-4122)

$ws.Range("H36").NumberFormat = "@"
$ws.Range("H36").Value = '66%'
$ws.Range("F36").Copy()
$ws.Range("H36").PasteSpecial(-4122)

$ws.Range("H37").NumberFormat = "@"
$ws.Range("H37").Value = '71%'
$ws.Range("F37").Copy()
$ws.Range("H37").PasteSpecial(-4122)

$ws.Range("H40").NumberFormat = "@"
$ws.Range("H40").Value = '74%'
$ws.Range("F40").Copy()
$ws.Range("H40").PasteSpecial(-4122)

$ws.Range("H43").NumberFormat = "@"
$ws.Range("H43").Value = '61%'
$ws.Range("F43").Copy()
$ws.Range("H43").PasteSpecial(-4122)

$ws.Range("H46").NumberFormat = "@"
$ws.Range("H46").Value = '45%'
$ws.Range("F46").Copy()
$ws.Range("H46").PasteSpecial(-4122)

$excel.CutCopyMode = $false
